# Update cryptocurrency price/volume figures to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'306.50"
$ws.Range("E2").Value = "'0.64%"
$ws.Range("D3").Value = "'36.73"
$ws.Range("E3").Value = "'2.68%"
$ws.Range("D4").Value = "'5.032"
$ws.Range("E4").Value = "'-1.08%"
$ws.Range("D5").Value = "'0.07852"
$ws.Range("E5").Value = "'0.12%"
$ws.Range("D6").Value = "'2.178"
$ws.Range("E6").Value = "'-3.58%"
$ws.Range("D7").Value = "'8.046"
$ws.Range("E7").Value = "'-0.97%"
$ws.Range("D8").Value = "'4.063"
$ws.Range("E8").Value = "'1.51%"
$ws.Range("D9").Value = "'0.9239"
$ws.Range("E9").Value = "'-0.31%"
$ws.Range("D10").Value = "'0.09954"
$ws.Range("E10").Value = "'1.01%"
$ws.Range("D11").Value = "'0.1875"
$ws.Range("E11").Value = "'2.85%"
$ws.Range("D12").Value = "'0.08689"
$ws.Range("E12").Value = "'-0.85%"
$ws.Range("D13").Value = "'0.03621"
$ws.Range("E13").Value = "'6.37%"
$ws.Range("D14").Value = "'0.09942"
$ws.Range("E14").Value = "'0.11%"
$ws.Range("D15").Value = "'0.001489"
$ws.Range("E15").Value = "'0.52%"
$ws.Range("D16").Value = "'0.005671"
$ws.Range("E16").Value = "'-1.85%"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.64%"
$ws.Range("D18").Value = "'2.335"
$ws.Range("E18").Value = "'9.75%"
$ws.Range("D19").Value = "'0.3451"
$ws.Range("E19").Value = "'0.56%"
$ws.Range("D20").Value = "'0.1348"
$ws.Range("E20").Value = "'2.10%"
$ws.Range("D21").Value = "'4.935"
$ws.Range("E21").Value = "'8.50%"
$ws.Range("D22").Value = "'0.2203"
$ws.Range("E22").Value = "'-1.43%"
$ws.Range("D23").Value = "'0.04626"
$ws.Range("E23").Value = "'-1.05%"
$ws.Range("D24").Value = "'0.005194"
$ws.Range("E24").Value = "'15.28%"
$ws.Range("D25").Value = "'0.001234"
$ws.Range("E25").Value = "'-0.45%"
$ws.Range("D26").Value = "'0.0001403"
$ws.Range("D27").Value = "'0.0002723"
$ws.Range("E27").Value = "'0.99%"
$ws.Range("D39").Value = "'0.01810"
$ws.Range("E39").Value = "'2.69%"
$ws.Range("D40").Value = "'0.04752"
$ws.Range("E40").Value = "'0.77%"
$ws.Range("D41").Value = "'0.007928"
$ws.Range("E41").Value = "'-1.10%"
$ws.Range("D42").Value = "'0.1410"
$ws.Range("E42").Value = "'-0.79%"
$ws.Range("D43").Value = "'0.007609"
$ws.Range("E43").Value = "'-10.84%"
$ws.Range("D44").Value = "'0.002184"
$ws.Range("E44").Value = "'-1.26%"
$ws.Range("E45").Value = "'10.68%"
$ws.Range("D46").Value = "'0.00006308"
$ws.Range("E46").Value = "'2.43%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.31%"
$ws.Range("D48").Value = "'0.0005813"
$ws.Range("E48").Value = "'0.22%"
$ws.Range("D49").Value = "'32.76"
$ws.Range("E49").Value = "'709.14%"
$ws.Range("E50").Value = "'0.20%"
$ws.Range("D51").Value = "'0.00002105"
$ws.Range("E51").Value = "'0.31%"
